# Add some assignment 6 grades
# Column layout: D=Assignment1, E=Assignment2, F=Assignment3, G=Assignment4,
# H=Assignment5, I=Assignment6, ... T=Total, U=Grade (percent of max row 2)
#
# Several students' "Assignment 4" score (column G) actually represented a
# combined 23 (or 22/18) that is being split out into the real Assignment 4
# (G, now 13 or 12), Assignment 5 (H, new = 10), and for a few rows
# Assignment 6 (I, new = 10 or 20) grades.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (max-possible-score row)
$ws.Cells.Item(2, 7).Value = 13   # G2
$ws.Cells.Item(2, 8).Value = 10   # H2 (new)
$ws.Cells.Item(2, 9).Value = 20   # I2 (new)

# Row 3
$ws.Cells.Item(3, 7).Value = 13   # G3
$ws.Cells.Item(3, 8).Value = 10   # H3 (new)

# Row 4
$ws.Cells.Item(4, 7).Value = 13   # G4
$ws.Cells.Item(4, 8).Value = 10   # H4 (new)

# Row 5
$ws.Cells.Item(5, 7).Value = 13   # G5
$ws.Cells.Item(5, 8).Value = 10   # H5 (new)
$ws.Cells.Item(5, 9).Value = 20   # I5 (new)

# Row 9
$ws.Cells.Item(9, 7).Value = 13   # G9
$ws.Cells.Item(9, 8).Value = 10   # H9 (new)

# Row 19
$ws.Cells.Item(19, 7).Value = 13  # G19
$ws.Cells.Item(19, 8).Value = 10  # H19 (new)

# Row 20
$ws.Cells.Item(20, 7).Value = 12  # G20
$ws.Cells.Item(20, 8).Value = 10  # H20 (new)

# Row 24
$ws.Cells.Item(24, 7).Value = 13  # G24
$ws.Cells.Item(24, 8).Value = 10  # H24 (new)
$ws.Cells.Item(24, 9).Value = 10  # I24 (new)

# Row 25
$ws.Cells.Item(25, 7).Value = 13  # G25
$ws.Cells.Item(25, 8).Value = 10  # H25 (new)

# Row 29
$ws.Cells.Item(29, 7).Value = 13  # G29
$ws.Cells.Item(29, 8).Value = 10  # H29 (new)

# Row 39
$ws.Cells.Item(39, 7).Value = 13  # G39
$ws.Cells.Item(39, 8).Value = 10  # H39 (new)

# Move the active selection to I18, matching the editor's final cursor spot.
$ws.Range("I18").Select() | Out-Null
